$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update summary figures ---
$ws.Range("E11").Value = 18170
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 2

# --- Reorder the two "MISHELLE" periods (row 18 = 2106, row 19 = 2105) ---
# so that after the rows above are removed, period 2105 ends up first
# (matching the new worker-record order), while keeping each row's own
# look (row 18 has the "normal" style, row 19 has the emphasized
# "last row" style).
$row18 = $ws.Range("B18:G18").Value()
$row19 = $ws.Range("B19:G19").Value()

$ws.Range("B18:G18").Value = $row19
$ws.Range("B19:G19").Value = $row18

# --- Remove the two obsolete worker rows (Rafael, Andres) ---
# Deleting these rows shifts the two "MISHELLE" rows up from 18/19 to 16/17,
# carrying their own row formatting with them.
$ws.Rows("16:17").Delete()
